# Fix the misspelled title "58 Interitance" -> "58 Inheritance" on the two
# slides that contain it (slide 2 and slide 3). The original title text was
# split across two runs ("58 " and "Interitance", the latter flagged with
# err="1" from the spell-checker); replacing the whole title text in one
# shot collapses it back into a single run using the first run's formatting
# (bold, yellow fill) and drops the stray spell-check "err" flag.

$p = $ppt.ActivePresentation

$slideIndexes = @(2, 3)

foreach ($idx in $slideIndexes) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item(1)
    $titleRange = $shape.TextFrame.TextRange
    $fullLen = $titleRange.Text.Length
    $titleRange.Characters(1, $fullLen).Text = "58 Inheritance"
}
